# Upgrade to PIO 3.8 beta5
#
# Adds a Fibonacci-style shared-formula column (B1:B10) to Sheet1:
#   B1 = 1
#   B2 = 2
#   B3:B10 = "=B1+B2" entered as one multi-cell formula so Excel treats it as
#            a relative formula that is re-based per row and serialised as a
#            single shared formula (t="shared") spanning B3:B10.
# Also leaves the B3:B10 range selected, matching the saved selection in the
# workbook's sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# Entering the same relative formula across the whole block makes Excel
# auto-adjust the references per row (B1+B2, B2+B3, ...) and store it as a
# single shared formula definition.
$ws.Range("B3:B10").Formula = "=B1+B2"

# Match the resulting selection/active cell recorded in the sheet view.
[void]$ws.Range("B3:B10").Select()
